# Equipment - Template.xlsx
#
# Commit: "Updated Hydrants to use new PM Shift Codes and added UserStatsu1."
# The ValidatedSource sheet gains a new "UserStatus1" header column, inserted
# right after "LifeCycleStatusCodeID" and before "ConditionRating" (i.e. the
# new column lands at AM1, pushing every column from the old AM onward one
# slot to the right, out to BI1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidatedSource")

# Insert a new blank column at AM (shifts AM:BH -> AN:BI).
$ws.Columns("AM:AM").Insert()

# Write the new header. The leading apostrophe marks the entry as explicit
# text (quote-prefixed) so the new cell picks up the same header formatting
# (bold font + text quote-prefix) already used by its sibling header cells.
$ws.Range("AM1").Value = "'UserStatus1"

# Reflect the author's post-edit cursor position on the sheet.
$ws.Range("AE1").Select()
